$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect to make the edits, then re-protect.
$ws.Unprotect()

# --- Update the confidential footer date (A16) ---
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# --- Update Weight (D) and Percent Change (E) columns for rows 2-13 ---
$ws.Range("D2").Value = 0.03065924183632839
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.02350282583683746
$ws.Range("E3").Value = 0.008654970760233915

$ws.Range("D4").Value = 0.05143565799490278
$ws.Range("E4").Value = 0.006958942240779198

$ws.Range("D5").Value = 0.136825924801251
$ws.Range("E5").Value = -0.006139925674584101

$ws.Range("D6").Value = 0.0317391377583176
$ws.Range("E6").Value = -0.001378359751895153

$ws.Range("D7").Value = 0.1198400228705861
$ws.Range("E7").Value = -0.0174786165860914

$ws.Range("D8").Value = 0.1028509626654526
$ws.Range("E8").Value = -0.001976639712488626

$ws.Range("D9").Value = 0.02939265680130898
$ws.Range("E9").Value = -0.003653338745686963

$ws.Range("D10").Value = 0.1263075841165464
$ws.Range("E10").Value = -0.003704389701796651

$ws.Range("D11").Value = 0.2446294127527848
$ws.Range("E11").Value = -0.02625601889706541

$ws.Range("D12").Value = 0.1028165725656839
$ws.Range("E12").Value = -0.01853975535168184

$ws.Range("D13").Value = 0.9999999999999999
$ws.Range("E13").Value = -0.01152489445514615

# Re-protect the sheet, matching the original protection options.
$ws.Protect("D382", $true, $true, $true)
